# Update NATMI LR-pair output (Apoe-Sorl1) with new TPM-derived values.
#
# The underlying source data changed for the "ECs" cluster:
#   - Apoe (ligand) average expression in ECs: 67.77251700000001 -> 47.23036199999999
#   - Sorl1 (receptor) average expression in ECs: 87.94215800000001 -> 0.033584
# All the dependent columns (totals, detection rates/cells, specificity scores,
# edge weights) are recomputed from those two base numbers. We write the fully
# recomputed values for every affected cell directly, matching the values
# produced by the upstream NATMI scripts with the refreshed TPM matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 47.23036199999999
$ws.Range("H2").Value = 141.691086
$ws.Range("I2").Value = 0.3244251370417807
$ws.Range("J2").Value = 0.3244251370417807
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.033584
$ws.Range("N2").Value = 0.100752
$ws.Range("O2").Value = 0.002172419590320632
$ws.Range("P2").Value = 0.002172419590320632
$ws.Range("Q2").Value = 1.586184477408
$ws.Range("R2").Value = 14.275660296672
$ws.Range("S2").Value = 0.0007047875233020201
$ws.Range("T2").Value = 0.0007047875233020202

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 47.23036199999999
$ws.Range("H3").Value = 141.691086
$ws.Range("I3").Value = 0.3244251370417807
$ws.Range("J3").Value = 0.3244251370417807
$ws.Range("O3").Value = 0.8824690642271135
$ws.Range("P3").Value = 0.8824690642271135
$ws.Range("Q3").Value = 644.3316648894779
$ws.Range("R3").Value = 5798.984984005301
$ws.Range("S3").Value = 0.2862951470970133
$ws.Range("T3").Value = 0.2862951470970133

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 47.23036199999999
$ws.Range("H4").Value = 141.691086
$ws.Range("I4").Value = 0.3244251370417807
$ws.Range("J4").Value = 0.3244251370417807
$ws.Range("O4").Value = 0.115358516182566
$ws.Range("P4").Value = 0.115358516182566
$ws.Range("Q4").Value = 84.22861242868798
$ws.Range("R4").Value = 758.0575118581919
$ws.Range("S4").Value = 0.03742520242146544
$ws.Range("T4").Value = 0.03742520242146544

# Row 5: FAPs -> ECs
$ws.Range("I5").Value = 0.4188548944674916
$ws.Range("J5").Value = 0.4188548944674916
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.033584
$ws.Range("N5").Value = 0.100752
$ws.Range("O5").Value = 0.002172419590320632
$ws.Range("P5").Value = 0.002172419590320632
$ws.Range("Q5").Value = 2.047871931098666
$ws.Range("R5").Value = 18.430847379888
$ws.Range("S5").Value = 0.0009099285782428596
$ws.Range("T5").Value = 0.0009099285782428597

# Row 6: FAPs -> FAPs
$ws.Range("I6").Value = 0.4188548944674916
$ws.Range("J6").Value = 0.4188548944674916
$ws.Range("O6").Value = 0.8824690642271135
$ws.Range("P6").Value = 0.8824690642271135
$ws.Range("S6").Value = 0.3696264867676737
$ws.Range("T6").Value = 0.3696264867676737

# Row 7: FAPs -> MuSCs
$ws.Range("I7").Value = 0.4188548944674916
$ws.Range("J7").Value = 0.4188548944674916
$ws.Range("O7").Value = 0.115358516182566
$ws.Range("P7").Value = 0.115358516182566
$ws.Range("S7").Value = 0.04831847912157508
$ws.Range("T7").Value = 0.04831847912157508

# Row 8: MuSCs -> ECs
$ws.Range("I8").Value = 0.2567199684907278
$ws.Range("J8").Value = 0.2567199684907277
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.033584
$ws.Range("N8").Value = 0.100752
$ws.Range("O8").Value = 0.002172419590320632
$ws.Range("P8").Value = 0.002172419590320632
$ws.Range("Q8").Value = 1.255159303541333
$ws.Range("R8").Value = 11.296433731872
$ws.Range("S8").Value = 0.0005577034887757523
$ws.Range("T8").Value = 0.0005577034887757523

# Row 9: MuSCs -> FAPs
$ws.Range("I9").Value = 0.2567199684907278
$ws.Range("J9").Value = 0.2567199684907277
$ws.Range("O9").Value = 0.8824690642271135
$ws.Range("P9").Value = 0.8824690642271135
$ws.Range("S9").Value = 0.2265474303624266
$ws.Range("T9").Value = 0.2265474303624265

# Row 10: MuSCs -> MuSCs
$ws.Range("I10").Value = 0.2567199684907278
$ws.Range("J10").Value = 0.2567199684907277
$ws.Range("O10").Value = 0.115358516182566
$ws.Range("P10").Value = 0.115358516182566
$ws.Range("S10").Value = 0.02961483463952544
$ws.Range("T10").Value = 0.02961483463952544
